# Atualizacao de bases das ligas - aplica swaps/alteracoes de linhas na aba Mexico Liga MX
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mexico Liga MX")


# Row 72
$ws.Range("B72").Value2 = 6754048
$ws.Range("E72").Value = "Atletico San Luis"
$ws.Range("F72").Value = "Mazatlan FC"
$ws.Range("G72").Value2 = 3
$ws.Range("I72").Value = "H"
$ws.Range("J72").Value2 = 1.615
$ws.Range("K72").Value2 = 4
$ws.Range("L72").Value2 = 4.5
$ws.Range("M72").Value2 = 1.6
$ws.Range("N72").Value2 = 4.5
$ws.Range("O72").Value2 = 5
$ws.Range("P72").Value2 = -1
$ws.Range("Q72").Value2 = 1.95
$ws.Range("R72").Value2 = 1.9
$ws.Range("S72").Value2 = 3
$ws.Range("T72").Value2 = 1.925
$ws.Range("U72").Value2 = 1.925
$ws.Range("V72").Value2 = 0.6000000000000001
$ws.Range("X72").Value2 = -1
$ws.Range("Y72").Value2 = 0
$ws.Range("Z72").Value2 = 0
$ws.Range("AA72").Value2 = 0.925

# Row 73
$ws.Range("B73").Value2 = 6754049
$ws.Range("E73").Value = "Juarez FC"
$ws.Range("F73").Value = "Atlas"
$ws.Range("G73").Value2 = 1
$ws.Range("I73").Value = "A"
$ws.Range("J73").Value2 = 2.75
$ws.Range("K73").Value2 = 3.25
$ws.Range("L73").Value2 = 2.375
$ws.Range("M73").Value2 = 2.6
$ws.Range("N73").Value2 = 3.2
$ws.Range("O73").Value2 = 2.8
$ws.Range("P73").Value2 = 0
$ws.Range("Q73").Value2 = 1.85
$ws.Range("R73").Value2 = 2
$ws.Range("S73").Value2 = 2.25
$ws.Range("T73").Value2 = 2.1
$ws.Range("U73").Value2 = 1.775
$ws.Range("V73").Value2 = -1
$ws.Range("X73").Value2 = 1.8
$ws.Range("Y73").Value2 = -1
$ws.Range("Z73").Value2 = 1
$ws.Range("AA73").Value2 = 1.1

# Row 75
$ws.Range("B75").Value2 = 6754051
$ws.Range("E75").Value = "Leon"
$ws.Range("F75").Value = "Tijuana"
$ws.Range("G75").Value2 = 1
$ws.Range("I75").Value = "H"
$ws.Range("J75").Value2 = 1.571
$ws.Range("K75").Value2 = 4
$ws.Range("L75").Value2 = 4.75
$ws.Range("M75").Value2 = 1.5
$ws.Range("N75").Value2 = 4.75
$ws.Range("O75").Value2 = 6
$ws.Range("P75").Value2 = -1.25
$ws.Range("Q75").Value2 = 2
$ws.Range("R75").Value2 = 1.85
$ws.Range("S75").Value2 = 3
$ws.Range("T75").Value2 = 1.925
$ws.Range("U75").Value2 = 1.925
$ws.Range("V75").Value2 = 0.5
$ws.Range("W75").Value2 = -1
$ws.Range("Z75").Value2 = 0.425
$ws.Range("AB75").Value2 = 0.925

# Row 76
$ws.Range("B76").Value2 = 6754052
$ws.Range("E76").Value = "Chivas Guadalajara"
$ws.Range("F76").Value = "Pachuca"
$ws.Range("G76").Value2 = 0
$ws.Range("I76").Value = "D"
$ws.Range("J76").Value2 = 2
$ws.Range("K76").Value2 = 3.3
$ws.Range("L76").Value2 = 3.4
$ws.Range("M76").Value2 = 2.2
$ws.Range("N76").Value2 = 3.2
$ws.Range("O76").Value2 = 3.5
$ws.Range("P76").Value2 = -0.25
$ws.Range("Q76").Value2 = 1.875
$ws.Range("R76").Value2 = 1.975
$ws.Range("S76").Value2 = 2.5
$ws.Range("T76").Value2 = 2.025
$ws.Range("U76").Value2 = 1.825
$ws.Range("V76").Value2 = -1
$ws.Range("W76").Value2 = 2.2
$ws.Range("Z76").Value2 = 0.4875
$ws.Range("AB76").Value2 = 0.825

# Row 92
$ws.Range("B92").Value2 = 6754065
$ws.Range("E92").Value = "Necaxa"
$ws.Range("F92").Value = "Cruz Azul"
$ws.Range("G92").Value2 = 1
$ws.Range("H92").Value2 = 3
$ws.Range("I92").Value = "A"
$ws.Range("J92").Value2 = 2.375
$ws.Range("K92").Value2 = 3.3
$ws.Range("L92").Value2 = 2.8
$ws.Range("M92").Value2 = 3.5
$ws.Range("O92").Value2 = 2.1
$ws.Range("P92").Value2 = 0.25
$ws.Range("Q92").Value2 = 2
$ws.Range("R92").Value2 = 1.85
$ws.Range("T92").Value2 = 1.9
$ws.Range("U92").Value2 = 1.95
$ws.Range("V92").Value2 = -1
$ws.Range("X92").Value2 = 1.1
$ws.Range("Y92").Value2 = -1
$ws.Range("Z92").Value2 = 0.8500000000000001
$ws.Range("AA92").Value2 = 0.8999999999999999

# Row 93
$ws.Range("B93").Value2 = 6754066
$ws.Range("E93").Value = "Unam Pumas"
$ws.Range("F93").Value = "Queretaro"
$ws.Range("G93").Value2 = 4
$ws.Range("H93").Value2 = 0
$ws.Range("I93").Value = "H"
$ws.Range("J93").Value2 = 1.727
$ws.Range("K93").Value2 = 3.5
$ws.Range("L93").Value2 = 4.5
$ws.Range("M93").Value2 = 1.8
$ws.Range("O93").Value2 = 4.5
$ws.Range("P93").Value2 = -0.75
$ws.Range("Q93").Value2 = 2.025
$ws.Range("R93").Value2 = 1.825
$ws.Range("T93").Value2 = 1.825
$ws.Range("U93").Value2 = 2.025
$ws.Range("V93").Value2 = 0.8
$ws.Range("X93").Value2 = -1
$ws.Range("Y93").Value2 = 1.025
$ws.Range("Z93").Value2 = -1
$ws.Range("AA93").Value2 = 0.825

# Row 303
$ws.Range("B303").Value2 = 7745553
$ws.Range("E303").Value = "Unam Pumas"
$ws.Range("F303").Value = "Leon"
$ws.Range("G303").Value2 = 1
$ws.Range("H303").Value2 = 0
$ws.Range("J303").Value2 = 2.2
$ws.Range("K303").Value2 = 3.5
$ws.Range("L303").Value2 = 3
$ws.Range("M303").Value2 = 1.909
$ws.Range("N303").Value2 = 3.8
$ws.Range("O303").Value2 = 3.8
$ws.Range("P303").Value2 = -0.5
$ws.Range("Q303").Value2 = 1.975
$ws.Range("R303").Value2 = 1.875
$ws.Range("S303").Value2 = 3
$ws.Range("T303").Value2 = 1.8
$ws.Range("U303").Value2 = 2.05
$ws.Range("V303").Value2 = 0.909
$ws.Range("Y303").Value2 = 0.9750000000000001
$ws.Range("AA303").Value2 = -1
$ws.Range("AB303").Value2 = 1.05

# Row 304
$ws.Range("B304").Value2 = 7745552
$ws.Range("E304").Value = "Atlas"
$ws.Range("F304").Value = "Atletico San Luis"
$ws.Range("G304").Value2 = 2
$ws.Range("H304").Value2 = 1
$ws.Range("J304").Value2 = 1.833
$ws.Range("K304").Value2 = 3.6
$ws.Range("L304").Value2 = 4.2
$ws.Range("M304").Value2 = 2.375
$ws.Range("N304").Value2 = 3.4
$ws.Range("O304").Value2 = 3
$ws.Range("P304").Value2 = -0.25
$ws.Range("Q304").Value2 = 2.05
$ws.Range("R304").Value2 = 1.8
$ws.Range("S304").Value2 = 2.75
$ws.Range("T304").Value2 = 1.85
$ws.Range("U304").Value2 = 2
$ws.Range("V304").Value2 = 1.375
$ws.Range("Y304").Value2 = 1.05
$ws.Range("AA304").Value2 = 0.425
$ws.Range("AB304").Value2 = -0.5

# Row 318
$ws.Range("B318").Value2 = 7612939
$ws.Range("E318").Value = "Toluca"
$ws.Range("F318").Value = "Cruz Azul"
$ws.Range("G318").Value2 = 0
$ws.Range("I318").Value = "A"
$ws.Range("J318").Value2 = 2
$ws.Range("K318").Value2 = 3.5
$ws.Range("L318").Value2 = 3.6
$ws.Range("M318").Value2 = 2.5
$ws.Range("N318").Value2 = 3.4
$ws.Range("O318").Value2 = 2.7
$ws.Range("P318").Value2 = 0
$ws.Range("S318").Value2 = 2.75
$ws.Range("T318").Value2 = 1.825
$ws.Range("U318").Value2 = 2.025
$ws.Range("V318").Value2 = -1
$ws.Range("X318").Value2 = 1.7
$ws.Range("Y318").Value2 = -1
$ws.Range("Z318").Value2 = 0.9750000000000001
$ws.Range("AA318").Value2 = -1
$ws.Range("AB318").Value2 = 1.025

# Row 319
$ws.Range("B319").Value2 = 7612937
$ws.Range("E319").Value = "Tigres UANL"
$ws.Range("F319").Value = "Tijuana"
$ws.Range("G319").Value2 = 4
$ws.Range("I319").Value = "H"
$ws.Range("J319").Value2 = 1.5
$ws.Range("K319").Value2 = 4.2
$ws.Range("L319").Value2 = 6.5
$ws.Range("M319").Value2 = 1.533
$ws.Range("N319").Value2 = 4.5
$ws.Range("O319").Value2 = 5.5
$ws.Range("P319").Value2 = -1
$ws.Range("S319").Value2 = 3
$ws.Range("T319").Value2 = 2.025
$ws.Range("U319").Value2 = 1.825
$ws.Range("V319").Value2 = 0.5329999999999999
$ws.Range("X319").Value2 = -1
$ws.Range("Y319").Value2 = 0.875
$ws.Range("Z319").Value2 = -1
$ws.Range("AA319").Value2 = 1.025
$ws.Range("AB319").Value2 = -1

# Row 322
$ws.Range("B322").Value2 = 7612941
$ws.Range("E322").Value = "Necaxa"
$ws.Range("F322").Value = "Monterrey"
$ws.Range("G322").Value2 = 2
$ws.Range("H322").Value2 = 5
$ws.Range("J322").Value2 = 3
$ws.Range("K322").Value2 = 3.5
$ws.Range("L322").Value2 = 2.25
$ws.Range("M322").Value2 = 3.25
$ws.Range("N322").Value2 = 3.4
$ws.Range("O322").Value2 = 2.2
$ws.Range("P322").Value2 = 0.25
$ws.Range("Q322").Value2 = 1.975
$ws.Range("R322").Value2 = 1.875
$ws.Range("S322").Value2 = 2.5
$ws.Range("T322").Value2 = 1.875
$ws.Range("U322").Value2 = 1.975
$ws.Range("X322").Value2 = 1.2
$ws.Range("Z322").Value2 = 0.875
$ws.Range("AA322").Value2 = 0.875
$ws.Range("AB322").Value2 = -1

# Row 323
$ws.Range("B323").Value2 = 8097226
$ws.Range("E323").Value = "Santos Laguna"
$ws.Range("F323").Value = "Atletico San Luis"
$ws.Range("G323").Value2 = 0
$ws.Range("H323").Value2 = 3
$ws.Range("J323").Value2 = 1.85
$ws.Range("K323").Value2 = 3.8
$ws.Range("L323").Value2 = 4
$ws.Range("M323").Value2 = 2.1
$ws.Range("N323").Value2 = 3.6
$ws.Range("O323").Value2 = 3.25
$ws.Range("P323").Value2 = -0.25
$ws.Range("Q323").Value2 = 1.825
$ws.Range("R323").Value2 = 2.025
$ws.Range("S323").Value2 = 2.75
$ws.Range("T323").Value2 = 2
$ws.Range("U323").Value2 = 1.85
$ws.Range("X323").Value2 = 2.25
$ws.Range("Z323").Value2 = 1.025
$ws.Range("AA323").Value2 = 0.5
$ws.Range("AB323").Value2 = -0.5

# Row 324
$ws.Range("M324").Value2 = 2.2
$ws.Range("Q324").Value2 = 1.9
$ws.Range("R324").Value2 = 1.95
$ws.Range("S324").Value2 = 2.5
$ws.Range("T324").Value2 = 1.95
$ws.Range("U324").Value2 = 1.9

# Row 325
$ws.Range("M325").Value2 = 2.375
$ws.Range("N325").Value2 = 3.4
$ws.Range("O325").Value2 = 2.9
$ws.Range("P325").Value2 = -0.25
$ws.Range("Q325").Value2 = 2.05
$ws.Range("R325").Value2 = 1.8
